# --- 🔄 Mise à jour du fichier Excel public ---
# Adds one new day (13-aug) as column BI to "Prix Spot", and one new row
# (2025-08-11) to "Gaz" and "CO2".

$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": add column BI (13-aug) with hourly prices ---
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Copy the header formatting (bold / border / centered) from the last
# existing header cell (BH1) onto the new one (BI1), then set its text.
$wsPrix.Range("BH1").Copy()
$wsPrix.Range("BI1").PasteSpecial(-4122)  # xlPasteFormats
$wsPrix.Range("BI1").Value = "13-aug"

$biValues = @{
    2  = 98.18000000000001
    3  = 90
    4  = 86.89
    5  = 84.70999999999999
    6  = 87.06999999999999
    7  = 92.76000000000001
    8  = 105
    9  = 103.24
    10 = 105.79
    11 = 94.09999999999999
    12 = 90
    13 = 71
    14 = 63.48
    15 = 47.01
    16 = 42.99
    17 = 70.90000000000001
    18 = 84.3
    19 = 96.64
    20 = 115.78
    21 = 150.09
    22 = 155.1
    23 = 153
    24 = 123.21
    25 = 106.27
}

foreach ($row in $biValues.Keys) {
    $wsPrix.Cells.Item($row, 61).Value = $biValues[$row]
}

# --- Sheet "Gaz": add row 58 (2025-08-11) ---
$wsGaz = $wb.Worksheets.Item("Gaz")

# Force the date to be kept as literal text (matching A2:A57), not
# auto-converted to a date serial, then drop the temporary text format so
# the cell ends up unstyled just like its neighbours.
$wsGaz.Range("A58").NumberFormat = "@"
$wsGaz.Range("A58").Value = "2025-08-11"
$wsGaz.Range("A58").ClearFormats()
$wsGaz.Range("B58").Value = 32.4

# --- Sheet "CO2": add row 58 (2025-08-11) ---
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("A58").NumberFormat = "@"
$wsCo2.Range("A58").Value = "2025-08-11"
$wsCo2.Range("A58").ClearFormats()
$wsCo2.Range("B58").Value = 71.73999999999999
